# "Generate Report for Handback" - record a handback transform failure for
# the 04d4478a-1aef-4457-b7ca-1864ccd6590d file across the Overview, zh-cn
# and de-de report sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status = "Handback transform failed"

# Overview sheet: row 3 is the 04d4478a-... record; zh-cn (E3) and de-de
# (F3) publish-status columns both flip from "Ready for handoff" to the
# failure status.
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status

# zh-cn detail sheet: row 3 Status column (C) + new Error Detail (P).
$zhcn.Range("C3").Value = $status
$zhcn.Range("P3").Value = "Handback file name: ju22u4yp.ggg is different with handoff file name: 04d4478a-1aef-4457-b7ca-1864ccd6590d.f02eecc8a8e7da9452f13c78f3da3a1a436b7fbf.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de detail sheet: row 3 Status column (C) + new Error Detail (P).
$dede.Range("C3").Value = $status
$dede.Range("P3").Value = "Handback file name: ju22u4yp.ggg is different with handoff file name: 04d4478a-1aef-4457-b7ca-1864ccd6590d.f02eecc8a8e7da9452f13c78f3da3a1a436b7fbf.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
